# Clash Identification and Resolving - PREUPDATING
# Update Sarah Obama's (MY004) Access Arrangement from "Alternative Site" to "N/A"
# and normalize the row heights of the data rows that no longer need to wrap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# --- Content change: Access Arrangement for MY004 / Sarah Obama (rows 17-21) ---
$ws.Range("F17:F21").Value = "N/A"

# --- Row height normalization (ht 30 -> 15.75) ---
$rowsToShrink = @(17,18,19,20,21,22,24,25,26,42,43,45,46,52,54,55,56)
foreach ($r in $rowsToShrink) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Update selection / view state to match the saved workbook state ---
$ws.Activate() | Out-Null
$ws.Range("I24").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
